# ---------------------------------------------------------------------------
# Applies the "added extra information on the new dataset ..." commit:
#   * adds a new worksheet "candidate list" (after "potential 5-8 angstrom")
#     with a small summary table of top clashing candidate structures
#   * tweaks a couple of window/selection bookkeeping bits on other sheets
#   * drops the stray codePage attribute on the 3rd text connection
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- 1. add the new "candidate list" worksheet as the LAST tab -------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "candidate list"

# --- 2. header row -----------------------------------------------------
$ws.Range("A1").Value = "PDB"
$ws.Range("B1").Value = "EMDB"
$ws.Range("C1").Value = "clashes"
$ws.Range("D1").Value = "chains"
$ws.Range("E1").Value = "unique chains"
$ws.Range("F1").Value = "res count"
$ws.Range("G1").Value = "resolution (Å)"
$ws.Range("H1").Value = "additional "
$ws.Range("I1").Value = "publication "

# --- 3. data rows (sorted by clash count, descending) -----------------
$ws.Range("A2").Value = "6r7i.pdb "
$ws.Range("B2").Value = "EMD-4742"
$ws.Range("C2").Value = 275
$ws.Range("D2").Value = 13
$ws.Range("E2").Value = 13
$ws.Range("F2").Value = 3911
$ws.Range("G2").Value = "5.9"
$ws.Range("H2").Value = "contains MSE residue and Zn ion"
$ws.Range("I2").Value = "doi: 10.1038/s41467-019-11772-y"
$ws.Range("J2").Value = "CSN–CRL2~N8 complex"

$ws.Range("A3").Value = "5n5z.pdb "
$ws.Range("B3").Value = "EMD-3591"
$ws.Range("C3").Value = 236
$ws.Range("D3").Value = 18
$ws.Range("E3").Value = 18
$ws.Range("F3").Value = 7778
$ws.Range("G3").Value = "7.7"
$ws.Range("I3").Value = "DOI: 10.1016/j.cell.2017.03.003"
$ws.Range("J3").Value = "RNA Polymerase I"

$ws.Range("A4").Value = "6n1q.pdb "
$ws.Range("B4").Value = "EMD-9317"
$ws.Range("C4").Value = 100
$ws.Range("D4").Value = 8
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 4088
$ws.Range("G4").Value = "5.16"
$ws.Range("I4").Value = "DOI: 10.7554/eLife.41215"
$ws.Range("J4").Value = "imers of gyrase A in complex with DNA illuminate"

$ws.Range("A5").Value = "6n8t.pdb "
$ws.Range("B5").Value = "EMD-0375"
$ws.Range("C5").Value = 86
$ws.Range("D5").Value = 6
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 5274
$ws.Range("G5").Value = "7.7"
$ws.Range("H5").Value = "contains ATP"
$ws.Range("I5").Value = "DOI: 10.1016/j.celrep.2018.12.037"
$ws.Range("J5").Value = "Hsp104 Protein Disaggregase"

$ws.Range("A6").Value = "5vj6.pdb "
$ws.Range("B6").Value = "EMD-8695"
$ws.Range("C6").Value = 77
$ws.Range("D6").Value = 18
$ws.Range("E6").Value = 6
$ws.Range("F6").Value = 3710
$ws.Range("G6").Value = "11.5"
$ws.Range("H6").Value = "Contains TYS residue"
$ws.Range("I6").Value = "DOI: 10.7554/eLife.27389"
$ws.Range("J6").Value = "HIV-1 Envelope trimer with antibodies"

$ws.Range("A7").Value = "6irf.pdb "
$ws.Range("B7").Value = "EMD-9715"
$ws.Range("C7").Value = 38
$ws.Range("D7").Value = 4
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 3376
$ws.Range("G7").Value = "5.1"
$ws.Range("I7").Value = "DOI: 10.1016/j.celrep.2018.11.071"
$ws.Range("J7").Value = "GluN1/GluN2A NMDA recepto"

$ws.Range("A8").Value = "6uc0.pdb "
$ws.Range("B8").Value = "EMD-20724"
$ws.Range("C8").Value = 27
$ws.Range("D8").Value = 8
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 2805
$ws.Range("G8").Value = "7.5"
$ws.Range("H8").Value = "contains ADP"
$ws.Range("I8").Value = "DOI: 10.2210/pdb6UC0/pdb"
$ws.Range("J8").Value = "S3D-cofilin bound to an actin filament"

$ws.Range("A9").Value = "5grs.pdb "
$ws.Range("B9").Value = "EMD-9537"
$ws.Range("C9").Value = 19
$ws.Range("D9").Value = 12
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 3084
$ws.Range("G9").Value = "5.4"
$ws.Range("I9").Value = "DOI: 10.1038/cr.2016.123"
$ws.Range("J9").Value = "SREBP-SCAP binding domains"

$ws.Range("A10").Value = "6c05.pdb "
$ws.Range("B10").Value = "EMD-7322"
$ws.Range("C10").Value = 15
$ws.Range("D10").Value = 7
$ws.Range("E10").Value = 6
$ws.Range("F10").Value = 3951
$ws.Range("G10").Value = "5.15"
$ws.Range("H10").Value = "Contains Zn and Mg"
$ws.Range("I10").Value = "DOI: 10.7554/eLife.34823"
$ws.Range("J10").Value = "Fdx/RbpA/σA-holo Complexes"

$ws.Range("A11").Value = "3j95.pdb "
$ws.Range("B11").Value = "EMD-6205"
$ws.Range("C11").Value = 13
$ws.Range("D11").Value = 6
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 4482
$ws.Range("G11").Value = "7.6"
$ws.Range("H11").Value = "contains ADP"
$ws.Range("I11").Value = "10.1038/nature14148"
$ws.Range("J11").Value = "ADP-bound N-ethylmaleimide sensitive factor"

# --- 4. a couple of font tweaks used on two of the cells above --------
$ws.Range("I2").Font.Name = "Arial"
$ws.Range("I2").Font.Size = 11
$ws.Range("I2").Font.Color = 0

$ws.Range("J6").Font.Name = "Calibri"
$ws.Range("J6").Font.Size = 12

# --- 5. column widths on the new sheet ---------------------------------
$ws.Range("D1").ColumnWidth = 10.83203125
$ws.Range("E1").ColumnWidth = 12.33203125
$ws.Range("F1").ColumnWidth = 10.83203125
$ws.Range("G1").ColumnWidth = 12.33203125
$ws.Range("H1").ColumnWidth = 38.5
$ws.Range("I1").ColumnWidth = 29.5

# --- 6. selection bookkeeping on the new sheet --------------------------
$ws.Range("D16").Select()

# --- 7. selection/active-sheet tweaks on pre-existing sheets -----------
$sheet2 = $wb.Worksheets.Item("Sheet2")
$sheet2.Range("D12").Select()

$sheet4 = $wb.Worksheets.Item("potential 5-8 angstrom")
$sheet4.Range("M22").Select()

# make the new sheet the active / selected tab
$ws.Select()

# --- 8. drop the stray codePage attribute on connection 3 ("clashes1") -
$conn = $wb.Connections.Item("clashes1")
$conn.ODBCConnection.CodePage = 0
